$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5513
$ws1.Range("F8").Value = 5429
$ws1.Range("F9").Value = 640
$ws1.Range("F10").Value = 11
$ws1.Range("F11").Value = 1399
$ws1.Range("F12").Value = 36

# Sheet "全部类型" (all types) - same events, same updated counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 5513
$ws4.Range("F9").Value = 5429
$ws4.Range("F10").Value = 640
$ws4.Range("F11").Value = 11
$ws4.Range("F12").Value = 1399
$ws4.Range("F13").Value = 36
